$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numara: (student number)
$ws.Range("M3").Value = 20215070019

# Ad Soyad: (student name)
$ws.Range("M4").Value = "KÜBRA ÇABUK"

# Bölüm: (department)
$ws.Range("M5").Value = "YBS"

# En küçük sayı (minimum of the number list)
$ws.Range("G6").Formula = "=MIN(D4:D17)"

# Move the active selection like the authored workbook
$ws.Range("M6").Select()
